$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the date formatting: "31-12-2024" -> "31/12/2024" (every cell using the old
#    text picks up the corrected shared string automatically).
$ws.Range("A2").Value = "31/12/2024"
$ws.Range("A3").Value = "31/12/2024"
$ws.Range("A4").Value = "31/12/2024"

# 2) Swap the two students so MAXIMILIANO JOAQUIN ALMONACID PEREZ (rut 21494146) is row 2
#    and FLAVIO ALEXANDER JARA LABRIN (rut 21075353) is row 3. Use Copy() (not Value=) so
#    the numeric-looking text (RUT/DV) keeps its original text type instead of being
#    re-parsed as a number.
$ws.Range("A2:G2").Copy($ws.Range("A20:G20"))
$ws.Range("A3:G3").Copy($ws.Range("A2:G2"))
$ws.Range("A20:G20").Copy($ws.Range("A3:G3"))
$ws.Range("A20:G20").Clear()

# 3) Remove the row for MATHIAS EDUARDO / DEUMACAN PULGAR (rut 21223313) entirely.
$ws.Rows.Item(4).Delete()

# 4) Add the thin border around every cell of the two remaining data rows.
$ws.Range("A2:G3").Borders.LineStyle = 1
